# Auto-generated Excel COM-interop script
# Applies the cell-level text/value updates for the cryptos list refresh
# described in the commit "Updated cryptos list on Sun Dec 24 10:33:27 UTC 2023
# with GitHub Actions".
#
# Notes:
#  - All target cells are plain text cells (t="inlineStr" in the original
#    workbook). Several Price values in column D look like numbers (e.g.
#    "109.71"); assigning such a string directly to Range.Value would make
#    Excel auto-convert it to a floating point number. To keep these as text
#    (matching the source data which also contains non-numeric looking
#    "thousands-dotted" values like "43.703.07" in the same column) we prefix
#    the handful of purely-numeric-looking values with a leading apostrophe,
#    which is the standard Excel convention for forcing text entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.703.07'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '2.286.65'
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("E4").Value = '  +0.41%  '
$ws.Range("D5").Value = '''109.71'
$ws.Range("E5").Value = '  +14.13%  '
$ws.Range("D6").Value = '''266.86'
$ws.Range("E6").Value = '  -0.54%  '
$ws.Range("E7").Value = '  +1.18%  '
$ws.Range("D9").Value = '''0.617'
$ws.Range("E9").Value = '  +1.02%  '
$ws.Range("D10").Value = '''47.34'
$ws.Range("E10").Value = '  +3.71%  '
$ws.Range("D11").Value = '''0.0945'
$ws.Range("E11").Value = '  +0.99%  '
$ws.Range("D12").Value = '''8.79'
$ws.Range("E12").Value = '  +10.67%  '
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("D14").Value = '''15.64'
$ws.Range("E14").Value = '  +1.93%  '
$ws.Range("D15").Value = '2.630.87'
$ws.Range("E15").Value = '  -0.18%  '
$ws.Range("D16").Value = '''0.843'
$ws.Range("E16").Value = '  -0.81%  '
$ws.Range("D17").Value = '2.289.40'
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").Value = '43.557.99'
$ws.Range("E18").Value = '  -0.12%  '
$ws.Range("D19").Value = '''0.0000109'
$ws.Range("E19").Value = '  +1.29%  '
$ws.Range("E20").Value = '  +6.27%  '
$ws.Range("D21").Value = '''72.26'
$ws.Range("E21").Value = '  -0.12%  '
$ws.Range("D22").Value = '''2.44'
$ws.Range("E22").Value = '  -4.68%  '
$ws.Range("D23").Value = '''231.86'
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("D24").Value = '''9.64'
$ws.Range("E24").Value = '  +5.88%  '
$ws.Range("D25").Value = '''2.77'
$ws.Range("E25").Value = '  +8.80%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").Value = '''11.56'
$ws.Range("E27").Value = '  +3.16%  '
$ws.Range("D28").Value = '''41.68'
$ws.Range("E28").Value = '  +4.17%  '
$ws.Range("D29").Value = '''3.40'
$ws.Range("E29").Value = '  -1.98%  '
$ws.Range("E30").Value = '  +1.52%  '
$ws.Range("D31").Value = '''176.04'
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("D32").Value = '''0.0925'
$ws.Range("E32").Value = '  +2.83%  '
$ws.Range("D33").Value = '''21.47'
$ws.Range("E33").Value = '  -2.42%  '
$ws.Range("D34").Value = '''5.59'
$ws.Range("E34").Value = '  +4.16%  '
$ws.Range("E35").Value = '  +0.80%  '
$ws.Range("D36").Value = '''4.69'
$ws.Range("E36").Value = '  +7.22%  '
$ws.Range("D37").Value = '''0.0360'
$ws.Range("E37").Value = '  +1.70%  '
$ws.Range("D38").Value = '''0.107'
$ws.Range("E38").Value = '  -0.90%  '
$ws.Range("E39").Value = '  +11.47%  '
$ws.Range("D40").Value = '''0.241'
$ws.Range("E40").Value = '  -0.18%  '
$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").Value = '''2.39'
$ws.Range("E41").Value = '  +3.57%  '
$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").Value = '''13.60'
$ws.Range("E42").Value = '  +10.59%  '
$ws.Range("D43").Value = '''70.74'
$ws.Range("E43").Value = '  +6.93%  '
$ws.Range("D44").Value = '''6.17'
$ws.Range("E44").Value = '  +18.74%  '
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D46").Value = '''1.39'
$ws.Range("E46").Value = '  +2.91%  '
$ws.Range("D47").Value = '''8.85'
$ws.Range("E47").Value = '  +0.61%  '
$ws.Range("E48").Value = '  -1.11%  '
$ws.Range("D49").Value = '''101.22'
$ws.Range("E49").Value = '  +3.98%  '
$ws.Range("E50").Value = '  +2.25%  '
$ws.Range("D51").Value = '''0.445'
$ws.Range("E51").Value = '  +5.82%  '
